$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values: revert to version 3
$ws.Range("A1").Value = 3
$ws.Range("B1").Value = 3
$ws.Range("C1").Value = 3

# Update the selected cell to E1 (this also updates the active cell / sqref)
$ws.Range("E1").Select()
